$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text so values like "299.51" or
# "0.3560" are not reinterpreted as numbers by Excel's auto-detection.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "23.478.51"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.646.73"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.44%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  +0.33%  "

# Row 6 - BNB
$ws.Range("D6").Value = "299.51"
$ws.Range("E6").Value = "  -1.50%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.31%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3560"
$ws.Range("E8").Value = "  -1.26%  "

# Row 9 - OKB
$ws.Range("D9").Value = "50.50"
$ws.Range("E9").Value = "  -2.89%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.08087"
$ws.Range("E10").Value = "  -1.36%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "1.218"
$ws.Range("E11").Value = "  -2.37%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13 - Solana
$ws.Range("D13").Value = "22.00"
$ws.Range("E13").Value = "  -2.17%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.401"
$ws.Range("E14").Value = "  -1.92%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "7.383"
$ws.Range("E15").Value = "  +0.05%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.00001200"
$ws.Range("E16").Value = "  -2.56%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.646.39"
$ws.Range("E17").Value = "  +0.52%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "97.07"
$ws.Range("E18").Value = "  +0.16%  "

# Row 19 - TRON
$ws.Range("D19").Value = "0.06981"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "6.759"
$ws.Range("E20").Value = "  +0.37%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "17.42"
$ws.Range("E21").Value = "  -1.00%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.30%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "12.50"
$ws.Range("E23").Value = "  -0.53%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "23.504.09"
$ws.Range("E24").Value = "  -0.57%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.480"
$ws.Range("E25").Value = "  -1.56%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "2.897"
$ws.Range("E26").Value = "  -7.01%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "20.93"
$ws.Range("E27").Value = "  -1.67%  "

# Row 28 - Monero
$ws.Range("D28").Value = "153.30"
$ws.Range("E28").Value = "  +0.66%  "

# Row 29 - HuobiToken
$ws.Range("D29").Value = "5.212"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "132.97"
$ws.Range("E30").Value = "  -1.35%  "

# Row 31 - WrappedliquidstakedEther2.0
$ws.Range("D31").Value = "1.834.13"
$ws.Range("E31").Value = "  +0.44%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "6.926"
$ws.Range("E32").Value = "  +2.14%  "

# Row 33 - WEMIXTOKEN
$ws.Range("D33").Value = "2.141"
$ws.Range("E33").Value = "  +4.42%  "

# Row 34 - FraxShare
$ws.Range("D34").Value = "11.86"
$ws.Range("E34").Value = "  +2.44%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "1.024"
$ws.Range("E35").Value = "  -5.88%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.02726"
$ws.Range("E36").Value = "  -2.85%  "

# Row 37 - Stellar
$ws.Range("D37").Value = "0.08735"
$ws.Range("E37").Value = "  -0.91%  "

# Row 38 - becomes Algorand (was Aptos)
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2448"
$ws.Range("E38").Value = "  -2.73%  "

# Row 39 - becomes Aptos (was Algorand)
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "13.38"
$ws.Range("E39").Value = "  +4.14%  "

# Row 40 - InternetComputer(DFINITY)
$ws.Range("D40").Value = "5.954"
$ws.Range("E40").Value = "  -2.30%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.06826"
$ws.Range("E41").Value = "  -2.96%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "0.6902"
$ws.Range("E42").Value = "  -2.29%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "1.318"
$ws.Range("E43").Value = "  -0.97%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "15.58"
$ws.Range("E44").Value = "  -2.33%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "0.6423"
$ws.Range("E45").Value = "  -1.33%  "

# Row 46 - Frax
$ws.Range("D46").Value = "0.9976"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47 - NEARProtocol
$ws.Range("D47").Value = "2.261"
$ws.Range("E47").Value = "  -3.31%  "

# Row 48 - PancakeSwap
$ws.Range("D48").Value = "3.925"
$ws.Range("E48").Value = "  -1.31%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "0.07792"
$ws.Range("E49").Value = "  -2.37%  "

# Row 50 - Quant
$ws.Range("D50").Value = "127.91"
$ws.Range("E50").Value = "  +0.10%  "

# Row 51 - Flow
$ws.Range("D51").Value = "1.165"
$ws.Range("E51").Value = "  -2.17%  "
